$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-27 Thursday", 2) | Out-Null
$d.Content.Find.Execute("67×38=2546", $true, $false, $false, $false, $false, $true, 1, $false, "23×93=2139", 2) | Out-Null
$d.Content.Find.Execute("98×66=6468", $true, $false, $false, $false, $false, $true, 1, $false, "54×64=3456", 2) | Out-Null
$d.Content.Find.Execute("43×63=2709", $true, $false, $false, $false, $false, $true, 1, $false, "59×63=3717", 2) | Out-Null
$d.Content.Find.Execute("95×48=4560", $true, $false, $false, $false, $false, $true, 1, $false, "57×68=3876", 2) | Out-Null
$d.Content.Find.Execute("12×74=888", $true, $false, $false, $false, $false, $true, 1, $false, "92×15=1380", 2) | Out-Null
$d.Content.Find.Execute("15×52=780", $true, $false, $false, $false, $false, $true, 1, $false, "69×53=3657", 2) | Out-Null
$d.Content.Find.Execute("51×32=1632", $true, $false, $false, $false, $false, $true, 1, $false, "46×91=4186", 2) | Out-Null
$d.Content.Find.Execute("15×51=765", $true, $false, $false, $false, $false, $true, 1, $false, "32×35=1120", 2) | Out-Null
$d.Content.Find.Execute("58×52=3016", $true, $false, $false, $false, $false, $true, 1, $false, "45×33=1485", 2) | Out-Null
$d.Content.Find.Execute("86×22=1892", $true, $false, $false, $false, $false, $true, 1, $false, "78×67=5226", 2) | Out-Null
$d.Content.Find.Execute("14×25=350", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=3726", 2) | Out-Null
$d.Content.Find.Execute("88×47=4136", $true, $false, $false, $false, $false, $true, 1, $false, "28×22=616", 2) | Out-Null
$d.Content.Find.Execute("23×48=1104", $true, $false, $false, $false, $false, $true, 1, $false, "77×40=3080", 2) | Out-Null
$d.Content.Find.Execute("44×50=2200", $true, $false, $false, $false, $false, $true, 1, $false, "28×41=1148", 2) | Out-Null
$d.Content.Find.Execute("57×79=4503", $true, $false, $false, $false, $false, $true, 1, $false, "31×78=2418", 2) | Out-Null
$d.Content.Find.Execute("59×68=4012", $true, $false, $false, $false, $false, $true, 1, $false, "57×85=4845", 2) | Out-Null
$d.Content.Find.Execute("61×47=2867", $true, $false, $false, $false, $false, $true, 1, $false, "19×68=1292", 2) | Out-Null
$d.Content.Find.Execute("92×19=1748", $true, $false, $false, $false, $false, $true, 1, $false, "66×59=3894", 2) | Out-Null
$d.Content.Find.Execute("33×22=726", $true, $false, $false, $false, $false, $true, 1, $false, "79×59=4661", 2) | Out-Null
$d.Content.Find.Execute("60×23=1380", $true, $false, $false, $false, $false, $true, 1, $false, "65×39=2535", 2) | Out-Null
$d.Content.Find.Execute("95×95=9025", $true, $false, $false, $false, $false, $true, 1, $false, "34×67=2278", 2) | Out-Null
$d.Content.Find.Execute("98×26=2548", $true, $false, $false, $false, $false, $true, 1, $false, "58×17=986", 2) | Out-Null
$d.Content.Find.Execute("96×74=7104", $true, $false, $false, $false, $false, $true, 1, $false, "92×52=4784", 2) | Out-Null
$d.Content.Find.Execute("73×59=4307", $true, $false, $false, $false, $false, $true, 1, $false, "33×83=2739", 2) | Out-Null
$d.Content.Find.Execute("58×51=2958", $true, $false, $false, $false, $false, $true, 1, $false, "81×70=5670", 2) | Out-Null
